$d = $word.ActiveDocument

# 1) Fix the double space between "Mails" and "and" in the first paragraph.
$d.Content.Find.Execute(
    "Daily morning Check Chats, Mails  and then Download latest build",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Daily morning Check Chats, Mails and then Download latest build", 2)

# 2) After the "22222 is a Valid Zip Code" paragraph (and the blank paragraph
#    that already follows it), add a new paragraph with the "Startup Items..."
#    text preceded by a new blank paragraph, keeping the document's original
#    trailing blank paragraph as the very last paragraph.
$trailingBlank = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailingBlank.Range.InsertBefore("Startup Items can also be enabled in Task manager`r")

$newTextPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newTextPara.Range.InsertParagraphBefore()
